{"js": "// Update two review-question answer paragraphs in the \"Regularization\"\n// section:\n//   1. \"L1 can lead co-efficients to 0; useful for feature selection.\"\n//      -> \"L1 can lead co-efficients to 0, useful for feature selection.\"\n//   2. \"Regularization is not limited regression. Whichever model is\n//       having co-efficients, regularization can be used. Even with\n//       decision tree.\"\n//      -> \"Regularization is not limited to regression. It can be used\n//       with any model that has co-efficients. Even with decision tree.\"\n\nconst body = context.document.body;\n\nconst edits = [\n  {\n    find: \"L1 can lead co-efficients to 0; useful for feature selection.\",\n    replace: \"L1 can lead co-efficients to 0, useful for feature selection.\"\n  },\n  {\n    find: \"Regularization is not limited regression. Whichever model is having co-efficients, regularization can be used. Even with decision tree.\",\n    replace: \"Regularization is not limited to regression. It can be used with any model that has co-efficients. Even with decision tree.\"\n  }\n];\n\nfor (const { find, replace } of edits) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find target text: \" + find);\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update two review-question answer paragraphs in the \"Regularization\"\n# section:\n#   1. \"L1 can lead co-efficients to 0; useful for feature selection.\"\n#      -> \"L1 can lead co-efficients to 0, useful for feature selection.\"\n#   2. \"Regularization is not limited regression. Whichever model is\n#       having co-efficients, regularization can be used. Even with\n#       decision tree.\"\n#      -> \"Regularization is not limited to regression. It can be used\n#       with any model that has co-efficients. Even with decision tree.\"\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceOne = 2\n\n$find1 = \"L1 can lead co-efficients to 0; useful for feature selection.\"\n$replace1 = \"L1 can lead co-efficients to 0, useful for feature selection.\"\n\n$r1 = $d.Content\n$r1.Find.ClearFormatting()\n$found1 = $r1.Find.Execute($find1, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replace1, $wdReplaceOne)\nif (-not $found1) {\n    throw \"Could not find target text: $find1\"\n}\n\n$find2 = \"Regularization is not limited regression. Whichever model is having co-efficients, regularization can be used. Even with decision tree.\"\n$replace2 = \"Regularization is not limited to regression. It can be used with any model that has co-efficients. Even with decision tree.\"\n\n$r2 = $d.Content\n$r2.Find.ClearFormatting()\n$found2 = $r2.Find.Execute($find2, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replace2, $wdReplaceOne)\nif (-not $found2) {\n    throw \"Could not find target text: $find2\"\n}\n"}
